$wb = $excel.ActiveWorkbook

# --- Sheet 1: Pediatric VFC Vaccine ---
$ws1 = $wb.Worksheets.Item(1)

# "DTaP" -> "DTaP/" (rows 2-5, shared text)
$ws1.Range("A2").Value = "DTaP/"
$ws1.Range("A3").Value = "DTaP/"
$ws1.Range("A4").Value = "DTaP/"
$ws1.Range("A5").Value = "DTaP/"

# "DTaP " -> "DTaP-Hib " (row 8, TriHIBit)
$ws1.Range("A8").Value = "DTaP-Hib "

# "Hepatitis B^" -> "Hepatitis B-Hib" (row 11, COMVAX)
$ws1.Range("A11").Value = "Hepatitis B-Hib"

# "Hepatitis A-Hepatitis B 18 only^" -> "Hepatitis A-Hepatitis B 18 only" (rows 16-17, Twinrix)
$ws1.Range("A16").Value = "Hepatitis A-Hepatitis B 18 only"
$ws1.Range("A17").Value = "Hepatitis A-Hepatitis B 18 only"

# "Hepatitis B PediatricAdolescent" -> "Hepatitis B Pediatric/Adolescent" (rows 18-21, ENGERIX B / RECOMBIVAX HB)
$ws1.Range("A18").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A19").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A20").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A21").Value = "Hepatitis B Pediatric/Adolescent"

# Split the combined ENGERIX B packaging text into 3 distinct rows (18,19,20)
$ws1.Range("D18").Value = "1 dose vials "
$ws1.Range("D19").Value = "10 pack - 1 dose vials "
$ws1.Range("D20").Value = "5 pack - 1 dose T-L syringes, No Needle "

# "Measles, Mumps and Rubella ( MMR)" -> "...( MMR)/" (row 28, MMRII)
$ws1.Range("A28").Value = "Measles, Mumps and Rubella ( MMR)/"

# "Pneumococcal7-valent (Pediatric)" -> "Pneumococcal 7-valent (Pediatric)" (row 29, Prevnar)
$ws1.Range("A29").Value = "Pneumococcal 7-valent (Pediatric)"

# "Tetanus  Diphtheria Toxoids^" -> "Tetanus  Diphtheria Toxoids" (rows 31-32, DECAVAC)
$ws1.Range("A31").Value = "Tetanus  Diphtheria Toxoids"
$ws1.Range("A32").Value = "Tetanus  Diphtheria Toxoids"

# Split the combined DECAVAC packaging text into 2 distinct rows (31,32)
$ws1.Range("D31").Value = "10 pack - 1 dose syringes No Needle "
$ws1.Range("D32").Value = "10 pack - 1 dose vials "

# "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis" -> same + "/" (rows 33-35, BOOSTRIX/ADACEL)
$ws1.Range("A33").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/"
$ws1.Range("A34").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/"
$ws1.Range("A35").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/"

# Split the combined BOOSTRIX packaging text into 2 distinct rows (33,34)
$ws1.Range("D33").Value = "10 pack - 1 dose vials "
$ws1.Range("D34").Value = "5 pack - 1 dose TL syringes, No Needle "

# --- Sheet 2: Adult Vaccine ---
$ws2 = $wb.Worksheets.Item(2)

# "Hepatitis A-Hepatitis B Adult^" -> "Hepatitis A-Hepatitis B Adult" (rows 6-7, Twinrix)
$ws2.Range("A6").Value = "Hepatitis A-Hepatitis B Adult"
$ws2.Range("A7").Value = "Hepatitis A-Hepatitis B Adult"

# "Tetanus  Diphtheria Toxoids^" -> "Tetanus  Diphtheria Toxoids" (row 13, shares the
# same text as sheet1 rows 31-32 / DECAVAC)
$ws2.Range("A13").Value = "Tetanus  Diphtheria Toxoids"

# --- Sheet 3: Influenza Vaccine ---
$ws3 = $wb.Worksheets.Item(3)

# Extra space fixes in the Influenza age-range labels
$ws3.Range("A2").Value = "Influenza   (Age 6 months and older)"
$ws3.Range("A3").Value = "Influenza  (Age 6-35 months)"
$ws3.Range("A4").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("A6").Value = "Influenza   (Age 4 years and older)"
$ws3.Range("A7").Value = "Influenza  (Age 18 years and older)"
$ws3.Range("A8").Value = "Influenza  Live, Intranasal (Age 5-49 years)"
